# Add data for 2025-03-10: update 2025 year-to-date (column L) totals
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1024
$ws.Range("L3").Value = 1040
$ws.Range("L4").Value = 286
$ws.Range("L5").Value = 69
$ws.Range("L6").Value = 1064
$ws.Range("L7").Value = 3483

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 61
$ws.Range("L3").Value = 66
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 9
$ws.Range("L7").Value = 111
$ws.Range("L8").Value = 210
$ws.Range("L11").Value = 55
$ws.Range("L16").Value = 8
$ws.Range("L19").Value = 107
$ws.Range("L20").Value = 89
$ws.Range("L23").Value = 37
$ws.Range("L29").Value = 161
$ws.Range("L31").Value = 38
$ws.Range("L32").Value = 6
$ws.Range("L33").Value = 151
$ws.Range("L34").Value = 20
$ws.Range("L37").Value = 119
$ws.Range("L41").Value = 14
$ws.Range("L48").Value = 56
$ws.Range("L49").Value = 20
$ws.Range("L52").Value = 65
$ws.Range("L54").Value = 74
$ws.Range("L55").Value = 39
$ws.Range("L63").Value = 15
$ws.Range("L65").Value = 71
$ws.Range("L67").Value = 127
$ws.Range("L68").Value = 10
$ws.Range("L76").Value = 41
$ws.Range("L78").Value = 54
$ws.Range("L79").Value = 102
$ws.Range("L80").Value = 14
$ws.Range("L83").Value = 77
$ws.Range("L84").Value = 35
$ws.Range("L85").Value = 181
$ws.Range("L86").Value = 26
$ws.Range("L95").Value = 50
$ws.Range("L96").Value = 31
$ws.Range("L97").Value = 40
$ws.Range("L98").Value = 29
$ws.Range("L99").Value = 55
$ws.Range("L101").Value = 3483

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 14
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 55
$ws.Range("L3").Value = 54
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 161

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 31
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 8
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 33
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 3
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 80
$ws.Range("L4").Value = 15
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 8
